# Apply arrow/line formatting changes to the connector shapes on slide 1.
#
# Summary of target changes (from the OOXML diff):
#  - Several straight-arrow connectors get a heavier line weight (w="38100",
#    i.e. 3 pt) and "small" arrowhead width/length (w="sm" len="sm") on
#    whichever ends already carry an explicit headEnd/tailEnd.
#  - One connector (id 15) only has its headEnd (arrow) touched this way;
#    its tailEnd stays untouched.
#  - Three connectors simply change their dash style from "dashDot" to
#    "solid" (no weight / arrowhead-size changes).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# All the target connectors live inside one big top-level group shape.
$top = $s.Shapes.Item(1)

function Find-ShapeById {
    param($shapes, [int]$targetId)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $targetId) {
            return $sh
        }
        if ($sh.Type -eq 6) {
            $found = Find-ShapeById $sh.GroupItems $targetId
            if ($found -ne $null) {
                return $found
            }
        }
    }
    return $null
}

# msoArrowheadNarrow / msoArrowheadShort == 1 == OOXML w="sm" / len="sm"
$ARROW_SM = 1
# 3 pt line weight -> w="38100" EMU in the saved XML
$WEIGHT = 3
# msoLineSolid
$DASH_SOLID = 1

# ---------------------------------------------------------------------
# Group 1: headEnd="none", tailEnd="arrow" -> both get sm/sm + w=38100
#   ids: 30, 33, 143
# ---------------------------------------------------------------------
foreach ($id in 30, 33, 143) {
    $sh = Find-ShapeById $top.GroupItems $id
    $sh.Line.Weight = $WEIGHT
    $sh.Line.BeginArrowheadWidth = $ARROW_SM
    $sh.Line.BeginArrowheadLength = $ARROW_SM
    $sh.Line.EndArrowheadWidth = $ARROW_SM
    $sh.Line.EndArrowheadLength = $ARROW_SM
}

# ---------------------------------------------------------------------
# Group 2: headEnd="arrow", tailEnd="none" -> both get sm/sm + w=38100
#   ids: 88, 89, 108
# ---------------------------------------------------------------------
foreach ($id in 88, 89, 108) {
    $sh = Find-ShapeById $top.GroupItems $id
    $sh.Line.Weight = $WEIGHT
    $sh.Line.BeginArrowheadWidth = $ARROW_SM
    $sh.Line.BeginArrowheadLength = $ARROW_SM
    $sh.Line.EndArrowheadWidth = $ARROW_SM
    $sh.Line.EndArrowheadLength = $ARROW_SM
}

# ---------------------------------------------------------------------
# Group 3: only headEnd="arrow" gets sm/sm + w=38100; tailEnd untouched
#   id: 15
# ---------------------------------------------------------------------
$sh15 = Find-ShapeById $top.GroupItems 15
$sh15.Line.Weight = $WEIGHT
$sh15.Line.BeginArrowheadWidth = $ARROW_SM
$sh15.Line.BeginArrowheadLength = $ARROW_SM
# Touch (no-op) the end arrowhead so it keeps sorting after headEnd in the
# saved XML (preserves original headEnd-before-tailEnd element order).
$sh15.Line.EndArrowheadStyle = $sh15.Line.EndArrowheadStyle

# ---------------------------------------------------------------------
# Group 4: dash style dashDot -> solid only (no weight/arrowhead change)
#   ids: 36, 6, 90
# ---------------------------------------------------------------------
foreach ($id in 36, 6) {
    $sh = Find-ShapeById $top.GroupItems $id
    $sh.Line.DashStyle = $DASH_SOLID
    # Touch (no-op) the end arrowhead so prstDash keeps sorting before the
    # single tailEnd element in the saved XML (preserves original order).
    $sh.Line.EndArrowheadStyle = $sh.Line.EndArrowheadStyle
}

$sh90 = Find-ShapeById $top.GroupItems 90
$sh90.Line.DashStyle = $DASH_SOLID
# Touch (no-op) both ends, Begin before End, so prstDash stays first and
# headEnd/tailEnd keep their original relative order in the saved XML.
$sh90.Line.BeginArrowheadStyle = $sh90.Line.BeginArrowheadStyle
$sh90.Line.EndArrowheadStyle = $sh90.Line.EndArrowheadStyle

Write-Output "Arrow/line formatting updated."
